$d = $word.ActiveDocument

# --- Part 1: Insert a new "Body" paragraph after the COVID-19 statistics paragraph ---
$covidPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*COVID-19 statistics and comparing the spread*") {
        $covidPara = $cand
        break
    }
}
if ($covidPara -eq $null) {
    throw "Could not find the COVID-19 statistics paragraph"
}
$r1 = $covidPara.Range
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3901D8FA" w14:textId="3A21B880" w:rsidR="00BF4F16" w:rsidRPr="00FA1791" w:rsidRDefault="00FA1791" w:rsidP="0023500F"><w:pPr><w:pStyle w:val="Body"/><w:rPr><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:iCs/></w:rPr><w:t>Here we will be looking at the COVID-19 statistics and comparing the spread of the virus in different areas and using the last census data to try to understand the various factors behind the spread of the virus.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:rPr><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:iCs/></w:rPr><w:t>To solve this problem, we have the COVID 19 case, death and Vaccine rates by UK region [1]. The ONS estimated age breakdown by region (as of August 2021). COVID 19 cases by age and region. A portion of the 2011 Census data</w:t></w:r><w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve"> showing the shared/unshared dwellings, number of cars, long term health and household deprivation. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Part 2: Replace the block of "Reference" styled paragraphs with the two new ones ---
$firstRef = $null
$lastRef = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t -like "*M.Ankerst, M.Breunig*") {
        $firstRef = $cand
    }
    if ($t -like "*N. Willems*") {
        $lastRef = $cand
    }
}
if ($firstRef -eq $null -or $lastRef -eq $null) {
    throw "Could not find reference paragraphs"
}
$r2 = $d.Range($firstRef.Range.Start, $lastRef.Range.End)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Reference"/></w:pPr><w:bookmarkStart w:id="0" w:name="_Ref371689630"/><w:r><w:t>Office for National Statistics; National Records of Scotland; Northern Ireland Statistics and Research Agency (2017): 2011 Census aggregate data. UK Data Service (Edition: February 2017). DOI: http://dx.doi.org/10.5257/census/aggregate-2011-2</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="Reference"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)

Write-Output "Edit applied successfully"
